$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 483.73334
$ws.Range("I33").Value = 535.75
$ws.Range("J33").Value = 424.2857
$ws.Range("K33").Value = 535.75
$ws.Range("L33").Value = 424.2857
$ws.Range("M33").Value = -306.75
$ws.Range("N33").Value = -882.2857
$ws.Range("H51").Value = 14862.647
$ws.Range("I51").Value = 12900
$ws.Range("J51").Value = 14985.3125
$ws.Range("K51").Value = 12900
$ws.Range("L51").Value = 14985.3125
$ws.Range("M51").Value = -12416
$ws.Range("N51").Value = -15953.3125
$ws.Range("H55").Value = 232.125
$ws.Range("I55").Value = 49
$ws.Range("J55").Value = 293.16666
$ws.Range("K55").Value = 49
$ws.Range("L55").Value = 293.16666
$ws.Range("M55").Value = 165
$ws.Range("N55").Value = -721.16666
$ws.Range("H80").Value = 1908.3158
$ws.Range("I80").Value = 215.3
$ws.Range("J80").Value = 3789.4443
$ws.Range("K80").Value = 645.9000000000001
$ws.Range("L80").Value = 11368.3329
$ws.Range("M80").Value = 352.0999999999999
$ws.Range("N80").Value = -13364.3329
$ws.Range("H83").Value = 1908.3158
$ws.Range("I83").Value = 215.3
$ws.Range("J83").Value = 3789.4443
$ws.Range("K83").Value = 1937.7
$ws.Range("L83").Value = 34104.9987
$ws.Range("M83").Value = 3054.3
$ws.Range("N83").Value = -44088.9987
$ws.Range("H100").Value = 4069.6667
$ws.Range("I100").Value = 4069.6667
$ws.Range("K100").Value = 4069.6667
$ws.Range("M100").Value = -3528.6667
$ws.Range("H107").Value = 931.2778
$ws.Range("I107").Value = 1071.7307
$ws.Range("J107").Value = 566.1
$ws.Range("K107").Value = 1071.7307
$ws.Range("L107").Value = 566.1
$ws.Range("M107").Value = 848.2692999999999
$ws.Range("N107").Value = -4406.1
$ws.Range("H131").Value = 2679.4
$ws.Range("I131").Value = 2732.6667
$ws.Range("J131").Value = 2599.5
$ws.Range("K131").Value = 8198.000100000001
$ws.Range("L131").Value = 7798.5
$ws.Range("M131").Value = -3158.000100000001
$ws.Range("N131").Value = -17878.5
$ws.Range("H137").Value = 8815.893
$ws.Range("I137").Value = 6004.857
$ws.Range("J137").Value = 17249
$ws.Range("K137").Value = 18014.571
$ws.Range("L137").Value = 51747
$ws.Range("M137").Value = -15464.571
$ws.Range("N137").Value = -56847
$ws.Range("H141").Value = 1393.9259
$ws.Range("I141").Value = 1235.2916
$ws.Range("K141").Value = 3705.8748
$ws.Range("M141").Value = 1474.1252

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3039.4
$ws.Range("I102").Value = 2487
$ws.Range("K102").Value = 2487
$ws.Range("M102").Value = -865
$ws.Range("H132").Value = 648849.0600000001
$ws.Range("I132").Value = 731103.9399999999
$ws.Range("J132").Value = 143569.14
$ws.Range("K132").Value = 2193311.82
$ws.Range("L132").Value = 430707.42
$ws.Range("M132").Value = -2190781.82
$ws.Range("N132").Value = -435767.42

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 16680392
$ws.Range("I80").Value = 2198.2856
$ws.Range("J80").Value = 25660958
$ws.Range("K80").Value = 2198.2856
$ws.Range("L80").Value = 25660958
$ws.Range("M80").Value = -1200.2856
$ws.Range("N80").Value = -25662954
$ws.Range("H81").Value = 21948.572
$ws.Range("J81").Value = 21948.572
$ws.Range("L81").Value = 21948.572
$ws.Range("N81").Value = -24070.572
$ws.Range("H83").Value = 16680392
$ws.Range("I83").Value = 2198.2856
$ws.Range("J83").Value = 25660958
$ws.Range("K83").Value = 10991.428
$ws.Range("L83").Value = 128304790
$ws.Range("M83").Value = -5999.428
$ws.Range("N83").Value = -128314774
$ws.Range("H84").Value = 21948.572
$ws.Range("J84").Value = 21948.572
$ws.Range("L84").Value = 65845.716
$ws.Range("N84").Value = -76453.716
$ws.Range("H105").Value = 2755.6
$ws.Range("I105").Value = 2788.3125
$ws.Range("J105").Value = 2624.75
$ws.Range("K105").Value = 2788.3125
$ws.Range("L105").Value = 2624.75
$ws.Range("M105").Value = -1041.3125
$ws.Range("N105").Value = -6118.75
$ws.Range("H107").Value = 12504876
$ws.Range("I107").Value = 14290715
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 14290715
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = -14288795
$ws.Range("N107").Value = -7840
$ws.Range("H134").Value = 882688.75
$ws.Range("I134").Value = 966024.2
$ws.Range("K134").Value = 2898072.6
$ws.Range("M134").Value = -2895537.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 24396010
$ws.Range("I58").Value = 32262362
$ws.Range("J58").Value = 10315.6
$ws.Range("K58").Value = 32262362
$ws.Range("L58").Value = 10315.6
$ws.Range("M58").Value = -32262159
$ws.Range("N58").Value = -10721.6
$ws.Range("H107").Value = 1638.2667
$ws.Range("I107").Value = 759
$ws.Range("J107").Value = 2407.625
$ws.Range("K107").Value = 759
$ws.Range("L107").Value = 2407.625
$ws.Range("M107").Value = 1161
$ws.Range("N107").Value = -6247.625
$ws.Range("H132").Value = 9858.352999999999
$ws.Range("I132").Value = 9698.833000000001
$ws.Range("J132").Value = 10241.2
$ws.Range("K132").Value = 29096.499
$ws.Range("L132").Value = 30723.6
$ws.Range("M132").Value = -26566.499
$ws.Range("N132").Value = -35783.60000000001
$ws.Range("H136").Value = 24396010
$ws.Range("I136").Value = 32262362
$ws.Range("J136").Value = 10315.6
$ws.Range("K136").Value = 96787086
$ws.Range("L136").Value = 30946.8
$ws.Range("M136").Value = -96784536
$ws.Range("N136").Value = -36046.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 747.38464
$ws.Range("I2").Value = 1146.2858
$ws.Range("J2").Value = 282
$ws.Range("K2").Value = 6877.714800000001
$ws.Range("L2").Value = 1692
$ws.Range("M2").Value = -6764.714800000001
$ws.Range("N2").Value = -1918
$ws.Range("H38").Value = 83.333336
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 83.333336
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 250.000008
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -944.000008
$ws.Range("H86").Value = 2785.5881
$ws.Range("I86").Value = 4962
$ws.Range("J86").Value = 2495.4
$ws.Range("K86").Value = 14886
$ws.Range("L86").Value = 7486.200000000001
$ws.Range("M86").Value = -13700
$ws.Range("N86").Value = -9858.200000000001
$ws.Range("H89").Value = 2785.5881
$ws.Range("I89").Value = 4962
$ws.Range("J89").Value = 2495.4
$ws.Range("K89").Value = 44658
$ws.Range("L89").Value = 22458.6
$ws.Range("M89").Value = -38730
$ws.Range("N89").Value = -34314.60000000001
$ws.Range("H98").Value = 934886.4399999999
$ws.Range("I98").Value = 1494083.2
$ws.Range("J98").Value = 2891.8333
$ws.Range("K98").Value = 4482249.6
$ws.Range("L98").Value = 8675.499899999999
$ws.Range("M98").Value = -4480751.6
$ws.Range("N98").Value = -11671.4999
$ws.Range("H122").Value = 222384.83
$ws.Range("I122").Value = 726
$ws.Range("J122").Value = 269049.84
$ws.Range("K122").Value = 6534
$ws.Range("L122").Value = 2421448.56
$ws.Range("M122").Value = -4084
$ws.Range("N122").Value = -2426348.56

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1318.3414
$ws.Range("I97").Value = 1277.129
$ws.Range("J97").Value = 1446.1
$ws.Range("K97").Value = 1277.129
$ws.Range("L97").Value = 1446.1
$ws.Range("M97").Value = -781.1289999999999
$ws.Range("N97").Value = -2438.1
$ws.Range("H126").Value = 30009268
$ws.Range("I126").Value = 38466336
$ws.Range("J126").Value = 14303285
$ws.Range("K126").Value = 115399008
$ws.Range("L126").Value = 42909855
$ws.Range("M126").Value = -115396538
$ws.Range("N126").Value = -42914795

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 7450
$ws.Range("I56").Value = 7450
$ws.Range("K56").Value = 7450
$ws.Range("M56").Value = -6759
$ws.Range("H61").Value = 6769.1665
$ws.Range("I61").Value = 6095.0557
$ws.Range("J61").Value = 8791.5
$ws.Range("K61").Value = 6095.0557
$ws.Range("L61").Value = 8791.5
$ws.Range("M61").Value = -5893.0557
$ws.Range("N61").Value = -9195.5
$ws.Range("H113").Value = 6769.1665
$ws.Range("I113").Value = 6095.0557
$ws.Range("J113").Value = 8791.5
$ws.Range("K113").Value = 6095.0557
$ws.Range("L113").Value = 8791.5
$ws.Range("M113").Value = -3925.0557
$ws.Range("N113").Value = -13131.5
$ws.Range("H132").Value = 2758.8708
$ws.Range("I132").Value = 2939.0625
$ws.Range("J132").Value = 2141.0715
$ws.Range("K132").Value = 8817.1875
$ws.Range("L132").Value = 6423.2145
$ws.Range("M132").Value = -6287.1875
$ws.Range("N132").Value = -11483.2145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3328.0454
$ws.Range("I126").Value = 2220.4443
$ws.Range("J126").Value = 8312.25
$ws.Range("K126").Value = 6661.3329
$ws.Range("L126").Value = 24936.75
$ws.Range("M126").Value = -4191.3329
$ws.Range("N126").Value = -29876.75
$ws.Range("H132").Value = 10018.546
$ws.Range("I132").Value = 8394.066000000001
$ws.Range("J132").Value = 13499.571
$ws.Range("K132").Value = 25182.198
$ws.Range("L132").Value = 40498.713
$ws.Range("M132").Value = -22652.198
$ws.Range("N132").Value = -45558.713

Write-Output "Applied 238 cell updates across 8 sheets"